$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "added Pi3 times to benchmark chart"
# Row 31 = Raspberry Pi 3 results for the 1/2/3 benchmark-count chart
# (chart4 / Sheet1!$B$31:$E$31). B31 (1 benchmark) already had a value;
# fill in the remaining 2/3/4-benchmark timings.
$ws.Range("C31").Value = 104.16
$ws.Range("D31").Value = 298
$ws.Range("E31").Value = 588.30999999999995

# Match number formatting of the rest of the row (0.00)
$ws.Range("E31").NumberFormat = "0.00"

# Scroll the sheet view down and leave the active selection on E32,
# matching the saved window state.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E32").Select()
